$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the product name typo: "Sam Brown Belt" -> "Sam Browne Belt"
# (rows 7 and 8, column A)
$ws.Range("A7").Value = "Sam Browne Belt"
$ws.Range("A8").Value = "Sam Browne Belt"

# Update the view: scroll so row 7 is at the top, and select D12
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D12").Select()
